$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $ok = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        Write-Host "WARNING: replacement not found for: $old"
    }
}

# 1. MVC Exchange — vice-versa (le Contrôleur demande au Modèle le taux, via son viewDidLoad()).
Replace-Text " (le Contrôleur demande au model le taux, via son " " (le Contrôleur demande au Modèle le taux, via son "

# 2. Le Contrôleur (ExchangeVC) va alors demander au Modèle (ExchangeService) le nouveau taux de change.
Replace-Text "Model (ExchangeService) le nouveau taux de change." "Modèle (ExchangeService) le nouveau taux de change."

# 3. Une fois le résultat obtenu, le Modèle l’envoie au Contrôleur via un callback.
Replace-Text "le résultat obtenu, le Model l’envoie" "le résultat obtenu, le Modèle l’envoie"

# 4. Full rewrite of the "Le contrôleur va alors..." sentence.
Replace-Text "Le contrôleur va alors se servir du taux de change précédemment obtenu et afficher le résultat du calcul à l’utilisateur." "Le Controller va alors demander au Modèle de réaliser le calcul, qui lui retournera le résultat pour qu’il soit affiché dans la Vue."

# 5. Le Contrôleur réceptionne l’Action et l’envoie au Modèle (TraductionService).
Replace-Text "Le Contrôleur réceptionne l’Action et l’envoie au Model (" "Le Contrôleur réceptionne l’Action et l’envoie au Modèle ("

# 6. Le Modèle va réaliser l’appel API ou il disposera...
Replace-Text "Le Model va réaliser l’appel API ou il" "Le Modèle va réaliser l’appel API ou il"

# 7. La traduction étant obtenue, le Modèle l’envoie...
Replace-Text "La traduction étant obtenue, le Model " "La traduction étant obtenue, le Modèle "

# 8. Le contrôleur demandera au Modèle (WeatherService) qui réalisera l’appel API...
Replace-Text "Le contrôleur demandera au Model (WeatherService) qui réalisera l’appel API afin d’obtenir la météo la plus récente." "Le contrôleur demandera au Modèle (WeatherService) qui réalisera l’appel API afin d’obtenir la météo la plus récente."

# 9. Une fois n’est pas coutume, le Modèle réalise l’appel API, s’appuie sur...
Replace-Text "Une fois n’est pas coutume, le Model réalise l’appel API, s’appuie sur " "Une fois n’est pas coutume, le Modèle réalise l’appel API, s’appuie sur "

# 10. Une fois les données reçues, le Modèle, grâce au callback,...
Replace-Text "Une fois les données reçues, le Model, grâce au callback," "Une fois les données reçues, le Modèle, grâce au callback,"

# 11. Cependant, ils restent assimilables à des Modèles et il est nécessaire...
Replace-Text "Cependant, ils restent assimilables à des Model " "Cependant, ils restent assimilables à des Modèles "

# 12. Via notre SessionFake, on lance l’appel API des Modèles (les fichiers Service).
Replace-Text "Via notre SessionFake, on lance l’appel API des Model " "Via notre SessionFake, on lance l’appel API des Modèles "
